$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 783.2857
$ws.Range("I92").Value = 555.82355
$ws.Range("J92").Value = 1750
$ws.Range("K92").Value = 555.82355
$ws.Range("L92").Value = 1750
$ws.Range("M92").Value = 692.17645
$ws.Range("N92").Value = -4246
$ws.Range("H98").Value = 1792.4615
$ws.Range("I98").Value = 1900.2
$ws.Range("J98").Value = 1433.3334
$ws.Range("K98").Value = 1900.2
$ws.Range("L98").Value = 1433.3334
$ws.Range("M98").Value = -402.2
$ws.Range("N98").Value = -4429.3334
$ws.Range("H122").Value = 1792.4615
$ws.Range("I122").Value = 1900.2
$ws.Range("J122").Value = 1433.3334
$ws.Range("K122").Value = 5700.6
$ws.Range("L122").Value = 4300.0002
$ws.Range("M122").Value = -3250.6
$ws.Range("N122").Value = -9200.0002
$ws.Range("H141").Value = 1625.2
$ws.Range("I141").Value = 1413.7693
$ws.Range("K141").Value = 4241.3079
$ws.Range("M141").Value = 938.6921000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3626.7407
$ws.Range("I32").Value = 3508.3774
$ws.Range("K32").Value = 3508.3774
$ws.Range("M32").Value = -3221.3774
$ws.Range("H101").Value = 41255.5
$ws.Range("J101").Value = 41255.5
$ws.Range("L101").Value = 41255.5
$ws.Range("N101").Value = -47745.5
$ws.Range("H110").Value = 13611.588
$ws.Range("I110").Value = 16616.416
$ws.Range("K110").Value = 16616.416
$ws.Range("M110").Value = -14571.416
$ws.Range("H132").Value = 4153.4062
$ws.Range("I132").Value = 4422.619
$ws.Range("J132").Value = 3639.4546
$ws.Range("K132").Value = 13267.857
$ws.Range("L132").Value = 10918.3638
$ws.Range("M132").Value = -10737.857
$ws.Range("N132").Value = -15978.3638
$ws.Range("H133").Value = 106444.125
$ws.Range("J133").Value = 106444.125
$ws.Range("L133").Value = 106444.125
$ws.Range("N133").Value = -111504.125
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 12274.777
$ws.Range("I26").Value = 12274.777
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 12274.777
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -11982.777
$ws.Range("N26").ClearContents()
$ws.Range("H81").Value = 23560.834
$ws.Range("J81").Value = 23560.834
$ws.Range("L81").Value = 23560.834
$ws.Range("N81").Value = -25682.834
$ws.Range("H84").Value = 23560.834
$ws.Range("J84").Value = 23560.834
$ws.Range("L84").Value = 70682.50199999999
$ws.Range("N84").Value = -81290.50199999999
$ws.Range("H103").Value = 28578.5
$ws.Range("J103").Value = 28578.5
$ws.Range("L103").Value = 28578.5
$ws.Range("N103").Value = -30922.5
$ws.Range("H132").Value = 105776.71
$ws.Range("J132").Value = 105776.71
$ws.Range("L132").Value = 105776.71
$ws.Range("N132").Value = -115896.71
$ws.Range("H140").Value = 96990
$ws.Range("J140").Value = 96990
$ws.Range("L140").Value = 96990
$ws.Range("N140").Value = -107350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 22600
$ws.Range("J43").Value = 22600
$ws.Range("L43").Value = 22600
$ws.Range("N43").Value = -22968
$ws.Range("H101").Value = 22600
$ws.Range("J101").Value = 22600
$ws.Range("L101").Value = 22600
$ws.Range("N101").Value = -29090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1290
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 1290
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3870
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -4340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 25993.8
$ws.Range("J39").Value = 25993.8
$ws.Range("L39").Value = 25993.8
$ws.Range("N39").Value = -27057.8
$ws.Range("H45").Value = 29833.334
$ws.Range("J45").Value = 29833.334
$ws.Range("L45").Value = 29833.334
$ws.Range("N45").Value = -30951.334
$ws.Range("H54").Value = 20178.4
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 20178.4
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 20178.4
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -20958.4
$ws.Range("H70").Value = 14380.917
$ws.Range("I70").Value = 4950.5
$ws.Range("K70").Value = 4950.5
$ws.Range("M70").Value = -4680.5
$ws.Range("H73").Value = 14380.917
$ws.Range("I73").Value = 4950.5
$ws.Range("K73").Value = 4950.5
$ws.Range("M73").Value = -4014.5
$ws.Range("H80").Value = 7066.1665
$ws.Range("J80").Value = 7066.1665
$ws.Range("L80").Value = 7066.1665
$ws.Range("N80").Value = -9062.166499999999
$ws.Range("H83").Value = 7066.1665
$ws.Range("J83").Value = 7066.1665
$ws.Range("L83").Value = 35330.8325
$ws.Range("N83").Value = -45314.8325
$ws.Range("H135").Value = 65194.25
$ws.Range("J135").Value = 65194.25
$ws.Range("L135").Value = 65194.25
$ws.Range("N135").Value = -75334.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9223
$ws.Range("I7").Value = 9912.4
$ws.Range("J7").Value = 7499.5
$ws.Range("K7").Value = 9912.4
$ws.Range("L7").Value = 7499.5
$ws.Range("M7").Value = -9800.4
$ws.Range("N7").Value = -7723.5
$ws.Range("H16").Value = 50001770
$ws.Range("I16").Value = 83334980
$ws.Range("K16").Value = 83334980
$ws.Range("M16").Value = -83334810
$ws.Range("H61").Value = 4300.7646
$ws.Range("I61").Value = 4623.2856
$ws.Range("J61").Value = 2795.6667
$ws.Range("K61").Value = 4623.2856
$ws.Range("L61").Value = 2795.6667
$ws.Range("M61").Value = -4421.2856
$ws.Range("N61").Value = -3199.6667
$ws.Range("H113").Value = 4300.7646
$ws.Range("I113").Value = 4623.2856
$ws.Range("J113").Value = 2795.6667
$ws.Range("K113").Value = 4623.2856
$ws.Range("L113").Value = 2795.6667
$ws.Range("M113").Value = -2453.2856
$ws.Range("N113").Value = -7135.6667
$ws.Range("H126").Value = 9223
$ws.Range("I126").Value = 9912.4
$ws.Range("J126").Value = 7499.5
$ws.Range("K126").Value = 29737.2
$ws.Range("L126").Value = 22498.5
$ws.Range("M126").Value = -27267.2
$ws.Range("N126").Value = -27438.5
$ws.Range("H136").Value = 56035.684
$ws.Range("I136").Value = 80575.46000000001
$ws.Range("K136").Value = 241726.38
$ws.Range("M136").Value = -239176.38

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 50389
$ws.Range("J93").Value = 50389
$ws.Range("L93").Value = 50389
$ws.Range("N93").Value = -55381
$ws.Range("H94").Value = 19500
$ws.Range("J94").Value = 19500
$ws.Range("L94").Value = 19500
$ws.Range("N94").Value = -21302
$ws.Range("H95").Value = 76625.63
$ws.Range("J95").Value = 76625.63
$ws.Range("L95").Value = 76625.63
$ws.Range("N95").Value = -82117.63
$ws.Range("H96").Value = 3431.3
$ws.Range("I96").Value = 3347.5715
$ws.Range("K96").Value = 3347.5715
$ws.Range("M96").Value = -1974.5715
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H101").Value = 10000
$ws.Range("J101").Value = 10000
$ws.Range("L101").Value = 10000
$ws.Range("N101").Value = -16490
$ws.Range("H103").Value = 14999.5
$ws.Range("J103").Value = 14999.5
$ws.Range("L103").Value = 14999.5
$ws.Range("N103").Value = -17343.5
$ws.Range("H105").Value = 37122.4
$ws.Range("J105").Value = 37122.4
$ws.Range("L105").Value = 37122.4
$ws.Range("N105").Value = -44110.4
$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676
$ws.Range("H125").Value = 29999
$ws.Range("J125").Value = 29999
$ws.Range("L125").Value = 29999
$ws.Range("N125").Value = -39839
$ws.Range("H137").Value = 91639.664
$ws.Range("J137").Value = 91639.664
$ws.Range("L137").Value = 91639.664
$ws.Range("N137").Value = -101839.664
$ws.Range("H140").Value = 69396.86
$ws.Range("J140").Value = 69396.86
$ws.Range("L140").Value = 69396.86
